$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $value) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue "D2" "26.840.30"

# Row 3 - Ethereum
Set-TextValue "D3" "1.815.11"
$ws.Range("E3").Value = "  -1.35%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.45%  "

# Row 5 - USDC
$ws.Range("E5").Value = "  -0.39%  "

# Row 6 - BNB
Set-TextValue "D6" "308.39"
$ws.Range("E6").Value = "  -2.05%  "

# Row 7 - XRP
Set-TextValue "D7" "0.4624"
$ws.Range("E7").Value = "  -2.57%  "

# Row 8 - Cardano
Set-TextValue "D8" "0.3640"
$ws.Range("E8").Value = "  -1.55%  "

# Row 9 - Dogecoin
Set-TextValue "D9" "0.07219"
$ws.Range("E9").Value = "  -3.28%  "

# Row 10 - Polygon
Set-TextValue "D10" "0.8563"
$ws.Range("E10").Value = "  -3.34%  "

# Row 11 - Solana
Set-TextValue "D11" "19.69"
$ws.Range("E11").Value = "  -3.76%  "

# Row 12 - TRON
Set-TextValue "D12" "0.07518"
$ws.Range("E12").Value = "  +2.33%  "

# Row 13 - WrappedEther
Set-TextValue "D13" "1.792.32"
$ws.Range("E13").Value = "  -6.63%  "

# Row 14 & 15 swap: Polkadot <-> Chainlink
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue "D14" "6.554"
$ws.Range("E14").Value = "  -0.53%  "

$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue "D15" "5.325"
$ws.Range("E15").Value = "  -2.38%  "

# Row 16 - Litecoin
Set-TextValue "D16" "91.72"
$ws.Range("E16").Value = "  -1.71%  "

# Row 17 - BinanceUSD
Set-TextValue "D17" "1.009"
$ws.Range("E17").Value = "  -0.15%  "

# Row 18 - ShibaInu
Set-TextValue "D18" "0.000008573"
$ws.Range("E18").Value = "  -2.89%  "

# Row 19 - Dai
Set-TextValue "D19" "1.009"
$ws.Range("E19").Value = "  -0.26%  "

# Row 20 - WrappedBTC
Set-TextValue "D20" "27.445.39"
$ws.Range("E20").Value = "  -0.20%  "

# Row 21 - Avalanche
Set-TextValue "D21" "14.40"
$ws.Range("E21").Value = "  -2.85%  "

# Row 22 - Uniswap
Set-TextValue "D22" "5.144"
$ws.Range("E22").Value = "  -3.42%  "

# Row 23 & 24 swap: Cosmos <-> WrappedliquidstakedEther2.0
$ws.Range("B23").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C23").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
Set-TextValue "D23" "2.168.98"
$ws.Range("E23").Value = "  +3.49%  "

$ws.Range("B24").Value = "Cosmos"
$ws.Range("C24").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue "D24" "10.50"
$ws.Range("E24").Value = "  -1.98%  "

# Row 25 - Monero
Set-TextValue "D25" "151.24"
$ws.Range("E25").Value = "  -0.74%  "

# Row 26 - Toncoin
$ws.Range("E26").Value = "  -2.87%  "

# Row 27 - EthereumClassic
Set-TextValue "D27" "18.07"
$ws.Range("E27").Value = "  -3.08%  "

# Row 28 - LidoDAOToken
Set-TextValue "D28" "2.066"
$ws.Range("E28").Value = "  -4.10%  "

# Row 29 - InternetComputer(DFINITY)
Set-TextValue "D29" "5.081"
$ws.Range("E29").Value = "  -3.43%  "

# Row 30 - BitcoinCash
Set-TextValue "D30" "114.91"
$ws.Range("E30").Value = "  -2.62%  "

# Row 31 - Stellar
Set-TextValue "D31" "0.08856"
$ws.Range("E31").Value = "  -1.67%  "

# Row 32 - HuobiToken
Set-TextValue "D32" "2.956"
$ws.Range("E32").Value = "  -0.01%  "

# Row 33 & 34 swap: Filecoin <-> ARBITRUM
$ws.Range("B33").Value = "ARBITRUM"
$ws.Range("C33").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue "D33" "1.130"
$ws.Range("E33").Value = "  -4.48%  "

$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue "D34" "4.401"
$ws.Range("E34").Value = "  -3.63%  "

# Row 35 - ImmutableX
Set-TextValue "D35" "0.7162"
$ws.Range("E35").Value = "  -5.47%  "

# Row 36 - Frax
$ws.Range("E36").Value = "  -0.55%  "

# Row 37 - TrustWalletToken
Set-TextValue "D37" "1.074"
$ws.Range("E37").Value = "  -2.88%  "

# Row 38 & 39 swap: RenderToken <-> Hedera
$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D38" "0.05233"
$ws.Range("E38").Value = "  -1.90%  "

$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D39" "2.424"
$ws.Range("E39").Value = "  +0.65%  "

# Row 41 - MXToken
Set-TextValue "D41" "2.916"
$ws.Range("E41").Value = "  -2.84%  "

# Row 42 - FraxShare
Set-TextValue "D42" "7.150"
$ws.Range("E42").Value = "  -2.62%  "

# Row 43 - TheSandbox
Set-TextValue "D43" "0.5126"
$ws.Range("E43").Value = "  -4.04%  "

# Row 44 - Algorand
Set-TextValue "D44" "0.1620"
$ws.Range("E44").Value = "  -2.52%  "

# Row 45 - Aptos
Set-TextValue "D45" "8.169"
$ws.Range("E45").Value = "  -4.17%  "

# Row 46 - Decentraland
Set-TextValue "D46" "0.4791"
$ws.Range("E46").Value = "  -2.71%  "

# Row 47 - PaxDollar
$ws.Range("E47").Value = "  -0.51%  "

# Row 48 & 49 swap: Quant <-> EnergySwap
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D48" "10.12"
$ws.Range("E48").Value = "  -4.67%  "

$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextValue "D49" "103.13"
$ws.Range("E49").Value = "  -1.93%  "

# Row 50 - NEARProtocol
$ws.Range("E50").Value = "  -3.97%  "

# Row 51 - Cronos
Set-TextValue "D51" "0.06216"
$ws.Range("E51").Value = "  -1.62%  "
